# Add a new "Season" column (AV) to the goalkeeper stats sheet, populated
# with 2023 for every data row (rows 2-16), matching the header row (row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Range("AV1").Value = "Season"

# Fill the Season value for every existing data row.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 48).Value = 2023
}

# Mirror the author's final selection/viewport from the saved workbook:
# a contiguous selection over the new column's data cells.
$null = $ws.Range("AV2:AV16").Select()
